$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-28 Monday" "2025-07-29 Tuesday"

Replace-Text "72÷5=" "59÷8="
Replace-Text "73÷3=" "45÷7="
Replace-Text "68÷7=" "34÷8="
Replace-Text "81÷8=" "38÷2="
Replace-Text "93÷2=" "50÷3="
Replace-Text "25÷8=" "45÷5="
Replace-Text "23÷2=" "47÷2="
Replace-Text "98÷6=" "56÷3="
Replace-Text "99÷3=" "10÷7="
Replace-Text "25÷2=" "32÷6="
Replace-Text "58÷6=" "26÷6="
Replace-Text "29÷9=" "29÷5="
Replace-Text "92÷6=" "44÷9="
Replace-Text "83÷7=" "66÷2="
Replace-Text "51÷8=" "14÷2="
Replace-Text "42÷5=" "35÷6="
Replace-Text "57÷8=" "79÷7="
Replace-Text "56÷2=" "95÷4="
Replace-Text "31÷3=" "77÷7="
Replace-Text "98÷5=" "38÷7="
Replace-Text "66÷3=" "41÷8="
Replace-Text "16÷3=" "90÷8="
Replace-Text "80÷2=" "90÷4="
Replace-Text "18÷9=" "43÷7="
Replace-Text "58÷5=" "23÷2="
